# New PO forecast model
# Applies updated weekly/monthly actuals and a refreshed PO forecast curve
# across the three sheets: "Weekly Quantity", "Monthly Trend", "PO Forecast".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws3 = $wb.Worksheets.Item("PO Forecast")

# --- "Weekly Quantity": append two new weekly actuals rows (39-40) ---------
$ws1.Cells.Item(39,1).Value = 45662.99999999999
$ws1.Cells.Item(39,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(39,2).Value = 10
$ws1.Cells.Item(40,1).Value = 45669.99999999999
$ws1.Cells.Item(40,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(40,2).Value = 10

# --- "Monthly Trend": append one new monthly actuals row (19) --------------
$ws2.Cells.Item(19,1).Value = 45688.99999999999
$ws2.Cells.Item(19,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(19,2).Value = 20

# --- "PO Forecast": refresh the forecast curve ------------------------------
# Existing rows 2-38 keep their dates; only the forecast quantity (col B)
# shifts to the new model's output. Rows 39-46 keep their positions but the
# model's date anchoring shifted by one week, so both columns move. Rows
# 47-48 are brand new weeks appended at the end of the forecast horizon.
$ws3.Cells.Item(2,2).Value = 76
$ws3.Cells.Item(3,2).Value = 75
$ws3.Cells.Item(4,2).Value = 74
$ws3.Cells.Item(5,2).Value = 74
$ws3.Cells.Item(6,2).Value = 74
$ws3.Cells.Item(7,2).Value = 74
$ws3.Cells.Item(8,2).Value = 74
$ws3.Cells.Item(9,2).Value = 73
$ws3.Cells.Item(10,2).Value = 73
$ws3.Cells.Item(11,2).Value = 72
$ws3.Cells.Item(12,2).Value = 72
$ws3.Cells.Item(13,2).Value = 72
$ws3.Cells.Item(14,2).Value = 72
$ws3.Cells.Item(20,2).Value = 70
$ws3.Cells.Item(21,2).Value = 70
$ws3.Cells.Item(22,2).Value = 66
$ws3.Cells.Item(23,2).Value = 66
$ws3.Cells.Item(24,2).Value = 66
$ws3.Cells.Item(25,2).Value = 66
$ws3.Cells.Item(26,2).Value = 65
$ws3.Cells.Item(27,2).Value = 64
$ws3.Cells.Item(28,2).Value = 64
$ws3.Cells.Item(29,2).Value = 64
$ws3.Cells.Item(30,2).Value = 63
$ws3.Cells.Item(31,2).Value = 63
$ws3.Cells.Item(32,2).Value = 62
$ws3.Cells.Item(33,2).Value = 62
$ws3.Cells.Item(34,2).Value = 62
$ws3.Cells.Item(35,2).Value = 62
$ws3.Cells.Item(36,2).Value = 62
$ws3.Cells.Item(37,2).Value = 61
$ws3.Cells.Item(38,2).Value = 61
$ws3.Cells.Item(39,1).Value = 45662.99999999999
$ws3.Cells.Item(39,2).Value = 61
$ws3.Cells.Item(40,1).Value = 45669.99999999999
$ws3.Cells.Item(40,2).Value = 60
$ws3.Cells.Item(41,1).Value = 45676.99999999999
$ws3.Cells.Item(41,2).Value = 60
$ws3.Cells.Item(42,1).Value = 45683.99999999999
$ws3.Cells.Item(42,2).Value = 60
$ws3.Cells.Item(43,1).Value = 45690.99999999999
$ws3.Cells.Item(43,2).Value = 60
$ws3.Cells.Item(44,1).Value = 45697.99999999999
$ws3.Cells.Item(44,2).Value = 60
$ws3.Cells.Item(45,1).Value = 45704.99999999999
$ws3.Cells.Item(45,2).Value = 59
$ws3.Cells.Item(46,1).Value = 45711.99999999999
$ws3.Cells.Item(46,2).Value = 59
$ws3.Cells.Item(47,1).Value = 45718.99999999999
$ws3.Cells.Item(47,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(47,2).Value = 59
$ws3.Cells.Item(48,1).Value = 45725.99999999999
$ws3.Cells.Item(48,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(48,2).Value = 59
